# Add a new QA test case row (Navigate-To-ToDoHistoryPage) to the "To Do App"
# test scenario table, mirroring the formatting of the row above it (row 4),
# then move the active selection to G4 as in the recorded session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is a pre-existing blank row right below the "Navigate To Do page" test
# case (row 4). Copy row 4's first five columns (A:E) into row 5 so the new
# row inherits the exact same cell formatting/styles used by row 4, then
# overwrite the copied values with the new test-case content.
$ws.Range("A4:E4").Copy($ws.Range("A5:E5"))

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "test_<Navigate-To-ToDoHistoryPage"
$ws.Range("C5").Value = "This is to test whether users are able to navigate to To Do History page"
$ws.Range("D5").Value = "NIL"
$ws.Range("E5").Value = "To Do History page is shown"

# Reflects the cell that was left selected in the saved workbook.
[void]$ws.Range("G4").Select()
